$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.433.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.092.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.23%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5218"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4372"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.14"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +16.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08868"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.154"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.55%  "
$ws.Range("E12").Value = "  -3.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.088.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.706"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.686"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "95.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.55%  "
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001121"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06592"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("E22").Value = "  -2.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.472.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("E24").Value = "  +1.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.338"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.330.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.570"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "131.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.53%  "
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.653"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.165"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.900"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02574"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06821"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.457"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2256"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6881"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.257"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6349"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.195"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.622"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.239"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.244"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.66%  "
